# Updated cryptos list with GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a cell with the default (unstyled) format as a style donor so that
# forcing text storage below does not leave a stray NumberFormat on the cell.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.351.22"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  +2.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.756.53"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  +3.20%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "115.60"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "330.76"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.55"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +1.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.29"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0829"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("E13").Value = "  +2.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.66"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +3.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.187.69"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  +2.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.744.27"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  +2.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.888"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "  +1.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.335.46"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  +2.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.61"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +3.25%  "

$ws.Range("E20").Value = "  +4.83%  "

$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.46"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +1.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.19"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  -3.19%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.91"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.35"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +2.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("E30").Value = "  -1.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.71"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.15"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  -0.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.65"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("E34").Value = "  +0.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.43"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  -1.24%  "

$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("E37").Value = "  +1.73%  "

$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("E39").Value = "  +2.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "129.36"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +3.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.74"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +4.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0353"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +10.37%  "

$ws.Range("E43").Value = "  +3.99%  "

$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("E45").Value = "  +3.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.117.47"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +9.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.26"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.55"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +3.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.09"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("E51").Value = "  +8.22%  "
